# Horarios actualizados Linea 141 - 161
# Applies the 14:53:07 scrape merge to the three schedule sheets.

$wb = $excel.ActiveWorkbook

$wsLP1912   = $wb.Worksheets.Item("LP1912")
$wsLP215    = $wb.Worksheets.Item("LP1912-215")
$wsL6203    = $wb.Worksheets.Item("6203-6173")

# ---------------------------------------------------------------------
# Sheet "LP1912": header refresh + re-sorted / newly merged data rows
# ---------------------------------------------------------------------
$wsLP1912.Range("A2").Value = "Última actualización: 14:53:07"
$wsLP1912.Range("A3").Value = "Total filas: 221"

$lp1912Rows = @(
    @{ Row = 128; A = "10:37:52"; B = "12:32"; C = "14_ABASTO";      D = 115 },
    @{ Row = 129; A = "11:47:17"; B = "12:32"; C = "23_HERNANDEZ";   D = 45  },

    @{ Row = 139; A = "11:34:59"; B = "12:47"; C = "14_ABASTO";      D = 73  },
    @{ Row = 140; A = "11:34:59"; B = "12:47"; C = "16_SANTA ANA";   D = 73  },
    @{ Row = 141; A = "11:34:59"; B = "12:47"; C = "15X38_ABASTO";   D = 73  },
    @{ Row = 142; A = "11:11:33"; B = "12:48"; C = "15X38_ABASTO";   D = 97  },
    @{ Row = 143; A = "11:47:17"; B = "12:48"; C = "14_ABASTO";      D = 61  },
    @{ Row = 144; A = "10:50:41"; B = "12:48"; C = "16_SANTA ANA";   D = 118 },

    @{ Row = 211; A = "14:53:07"; B = "15:53"; C = "10_OLMOS";       D = 60  },
    @{ Row = 212; A = "14:32:38"; B = "15:55"; C = "27_EL RETIRO";   D = 83  },
    @{ Row = 213; A = "14:12:26"; B = "15:56"; C = "27_EL RETIRO";   D = 104 },
    @{ Row = 214; A = "14:53:07"; B = "15:57"; C = "27_EL RETIRO";   D = 64  },
    @{ Row = 215; A = "14:12:26"; B = "16:05"; C = "14_ABASTO";      D = 113 },
    @{ Row = 216; A = "14:53:07"; B = "16:06"; C = "14_ABASTO";      D = 73  },
    @{ Row = 217; A = "14:53:07"; B = "16:13"; C = "16_SANTA ANA";   D = 80  },
    @{ Row = 218; A = "14:32:38"; B = "16:14"; C = "17_ROMERO";      D = 102 },
    @{ Row = 219; A = "14:32:38"; B = "16:17"; C = "10_OLMOS";       D = 105 },
    @{ Row = 220; A = "14:53:07"; B = "16:20"; C = "23_HERNANDEZ";   D = 87  },
    @{ Row = 221; A = "14:32:38"; B = "16:21"; C = "23_HERNANDEZ";   D = 109 },
    @{ Row = 222; A = "14:45:56"; B = "16:33"; C = "83_ALUAR";       D = 108 },
    @{ Row = 223; A = "14:53:07"; B = "16:34"; C = "83_ALUAR";       D = 101 },
    @{ Row = 224; A = "14:45:56"; B = "16:40"; C = "225_GOMEZ";      D = 115 },
    @{ Row = 225; A = "14:53:07"; B = "16:41"; C = "225_GOMEZ";      D = 108 },
    @{ Row = 226; A = "14:53:07"; B = "16:46"; C = "17_ROMERO";      D = 113 }
)

foreach ($r in $lp1912Rows) {
    $wsLP1912.Cells.Item($r.Row, 1).Value = $r.A
    $wsLP1912.Cells.Item($r.Row, 2).Value = $r.B
    $wsLP1912.Cells.Item($r.Row, 3).Value = $r.C
    $wsLP1912.Cells.Item($r.Row, 4).Value = $r.D
    $wsLP1912.Cells.Item($r.Row, 5).Value = "LP1912"
}

# ---------------------------------------------------------------------
# Sheet "LP1912-215": header timestamp refresh only
# ---------------------------------------------------------------------
$wsLP215.Range("A2").Value = "Última actualización: 14:53:07"

# ---------------------------------------------------------------------
# Sheet "6203-6173": header refresh + one new merged row
# ---------------------------------------------------------------------
$wsL6203.Range("A2").Value = "Última actualización: 14:53:07"
$wsL6203.Range("A3").Value = "Total filas: 33"

$wsL6203.Cells.Item(38, 1).Value = "14:53:07"
$wsL6203.Cells.Item(38, 2).Value = "16:30"
$wsL6203.Cells.Item(38, 3).Value = "215B_LP-P MOR-40 Y 115"
$wsL6203.Cells.Item(38, 4).Value = 97
$wsL6203.Cells.Item(38, 5).Value = "L6173"
